# Updated cryptos list with GitHub Actions: refresh Price (D) / Volume(1h) (E)
# columns for the coinranking.com snapshot on the sheet's only worksheet.
#
# Price cells are forced to Text (NumberFormat "@") before the assignment so
# that values such as "276.94" or "1.000" are stored verbatim instead of
# being auto-coerced into numbers (which would also silently drop
# significant trailing zeros, e.g. "0.5080" -> 0.508). The style is then
# reset back to "Normal" so no stray cell-format index is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  -3.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.813.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '276.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -8.02%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5080'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.69%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3513'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.35'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06674'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.00%  '
$ws.Range("E11").Value = '  -7.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8304'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07909'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.805.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.082'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008031'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '25.822.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.725'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.994'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.073'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.166'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.671'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '109.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.324'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.234'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08825'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04863'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7339'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.134'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.16%  '
$ws.Range("E36").Value = '  -3.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.154'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9995'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5214'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -12.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01841'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.301'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9576'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '113.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.182'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.069'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9995'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4573'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.36%  '
$ws.Range("E48").Value = '  -8.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.385'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.502'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.34%  '
